$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows (in descending order so row indices above don't shift)
# Original rows to remove (by their "even_MAG-GUT*.fa" content):
#   row 16 -> even_MAG-GUT6929.fa
#   row 15 -> even_MAG-GUT66949.fa
#   row 6  -> even_MAG-GUT12678.fa
#   row 5  -> even_MAG-GUT12269.fa
$ws.Rows("16:16").Delete()
$ws.Rows("15:15").Delete()
$ws.Rows("6:6").Delete()
$ws.Rows("5:5").Delete()
